# Auto-generated edit script applying the diff to Brynhildr_Profits workbook
# Updates LeveProfit/average-price columns (H-N) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 98997.5
$ws.Range("J57").Value = 98997.5
$ws.Range("L57").Value = 296992.5
$ws.Range("N57").Value = -297990.5
$ws.Range("H116").Value = 16820.047
$ws.Range("I116").Value = 4411.1816
$ws.Range("J116").Value = 30469.8
$ws.Range("K116").Value = 4411.1816
$ws.Range("L116").Value = 30469.8
$ws.Range("M116").Value = -969.1815999999999
$ws.Range("N116").Value = -37353.8
$ws.Range("H125").Value = 14347
$ws.Range("I125").Value = 2333
$ws.Range("J125").Value = 21555.4
$ws.Range("K125").Value = 20997
$ws.Range("L125").Value = 193998.6
$ws.Range("M125").Value = -18537
$ws.Range("N125").Value = -198918.6
$ws.Range("H132").Value = 7772.9443
$ws.Range("I132").Value = 8131.091
$ws.Range("K132").Value = 24393.273
$ws.Range("M132").Value = -21863.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 419.2857
$ws.Range("I4").Value = 265.27274
$ws.Range("J4").Value = 984
$ws.Range("K4").Value = 265.27274
$ws.Range("L4").Value = 984
$ws.Range("M4").Value = -149.27274
$ws.Range("N4").Value = -1216
$ws.Range("H45").Value = 1525.4348
$ws.Range("I45").Value = 1405.0667
$ws.Range("K45").Value = 1405.0667
$ws.Range("M45").Value = -1028.0667
$ws.Range("H74").Value = 963052.9399999999
$ws.Range("I74").Value = 1159419.2
$ws.Range("K74").Value = 1159419.2
$ws.Range("M74").Value = -1158545.2
$ws.Range("H77").Value = 963052.9399999999
$ws.Range("I77").Value = 1159419.2
$ws.Range("K77").Value = 5797096
$ws.Range("M77").Value = -5792728
$ws.Range("H122").Value = 25001496
$ws.Range("I122").Value = 33334542
$ws.Range("J122").Value = 2358.6
$ws.Range("K122").Value = 100003626
$ws.Range("L122").Value = 7075.799999999999
$ws.Range("M122").Value = -100001176
$ws.Range("N122").Value = -11975.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 3590
$ws.Range("I5").Value = 4425
$ws.Range("K5").Value = 4425
$ws.Range("M5").Value = -4312
$ws.Range("H7").Value = 417549.5
$ws.Range("I7").Value = 500059.6
$ws.Range("J7").Value = 4999
$ws.Range("K7").Value = 500059.6
$ws.Range("L7").Value = 4999
$ws.Range("M7").Value = -499946.6
$ws.Range("N7").Value = -5225
$ws.Range("H58").Value = 101500
$ws.Range("J58").Value = 101500
$ws.Range("L58").Value = 101500
$ws.Range("N58").Value = -102088
$ws.Range("H86").Value = 1349.7391
$ws.Range("I86").Value = 1402.7
$ws.Range("J86").Value = 996.6667
$ws.Range("K86").Value = 1402.7
$ws.Range("L86").Value = 996.6667
$ws.Range("M86").Value = -279.7
$ws.Range("N86").Value = -3242.6667
$ws.Range("H89").Value = 1349.7391
$ws.Range("I89").Value = 1402.7
$ws.Range("J89").Value = 996.6667
$ws.Range("K89").Value = 7013.5
$ws.Range("L89").Value = 4983.3335
$ws.Range("M89").Value = -1397.5
$ws.Range("N89").Value = -16215.3335
$ws.Range("H107").Value = 9616076
$ws.Range("I107").Value = 10417299
$ws.Range("K107").Value = 10417299
$ws.Range("M107").Value = -10415379
$ws.Range("H134").Value = 3408096.5
$ws.Range("I134").Value = 4818.512
$ws.Range("K134").Value = 14455.536
$ws.Range("M134").Value = -11920.536
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 79998.5
$ws.Range("J137").Value = 79998.5
$ws.Range("L137").Value = 79998.5
$ws.Range("N137").Value = -90198.5
$ws.Range("H138").Value = 89890
$ws.Range("J138").Value = 89890
$ws.Range("L138").Value = 89890
$ws.Range("N138").Value = -100170
$ws.Range("H139").Value = 40000
$ws.Range("I139").Value = 40000
$ws.Range("K139").Value = 40000
$ws.Range("M139").Value = -34860
$ws.Range("H141").Value = 299999
$ws.Range("J141").Value = 299999
$ws.Range("L141").Value = 299999
$ws.Range("N141").Value = -310359

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 60726.59
$ws.Range("I16").Value = 1989.7778
$ws.Range("J16").Value = 126805.5
$ws.Range("K16").Value = 1989.7778
$ws.Range("L16").Value = 126805.5
$ws.Range("M16").Value = -1702.7778
$ws.Range("N16").Value = -127379.5
$ws.Range("H22").Value = 819.95
$ws.Range("I22").Value = 395.375
$ws.Range("K22").Value = 395.375
$ws.Range("M22").Value = -45.375
$ws.Range("H42").Value = 15000
$ws.Range("I42").Value = 15000
$ws.Range("K42").Value = 15000
$ws.Range("M42").Value = -14407
$ws.Range("H62").Value = 3462.2
$ws.Range("J62").Value = 3612.1667
$ws.Range("L62").Value = 3612.1667
$ws.Range("N62").Value = -4860.1667
$ws.Range("H65").Value = 3462.2
$ws.Range("J65").Value = 3612.1667
$ws.Range("L65").Value = 18060.8335
$ws.Range("N65").Value = -24300.8335
$ws.Range("H113").Value = 60726.59
$ws.Range("I113").Value = 1989.7778
$ws.Range("J113").Value = 126805.5
$ws.Range("K113").Value = 1989.7778
$ws.Range("L113").Value = 126805.5
$ws.Range("M113").Value = 180.2221999999999
$ws.Range("N113").Value = -131145.5
$ws.Range("H132").Value = 2679.0344
$ws.Range("I132").Value = 2757.5789
$ws.Range("K132").Value = 8272.736699999999
$ws.Range("M132").Value = -5742.736699999999
$ws.Range("H133").Value = 67000
$ws.Range("J133").Value = 67000
$ws.Range("L133").Value = 67000
$ws.Range("N133").Value = -72060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 6514.143
$ws.Range("J68").Value = 8459.799999999999
$ws.Range("L68").Value = 25379.4
$ws.Range("N68").Value = -27001.4
$ws.Range("H70").Value = 668
$ws.Range("I70").Value = 668
$ws.Range("K70").Value = 2004
$ws.Range("M70").Value = -1689
$ws.Range("H71").Value = 6514.143
$ws.Range("J71").Value = 8459.799999999999
$ws.Range("L71").Value = 76138.2
$ws.Range("N71").Value = -84250.2
$ws.Range("H73").Value = 668
$ws.Range("I73").Value = 668
$ws.Range("K73").Value = 2004
$ws.Range("M73").Value = -912

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 10888.333
$ws.Range("J12").Value = 7666
$ws.Range("L12").Value = 7666
$ws.Range("N12").Value = -7946
$ws.Range("H70").Value = 25275.941
$ws.Range("I70").Value = 27138.3
$ws.Range("J70").Value = 22615.428
$ws.Range("K70").Value = 27138.3
$ws.Range("L70").Value = 22615.428
$ws.Range("M70").Value = -26868.3
$ws.Range("N70").Value = -23155.428
$ws.Range("H73").Value = 25275.941
$ws.Range("I73").Value = 27138.3
$ws.Range("J73").Value = 22615.428
$ws.Range("K73").Value = 27138.3
$ws.Range("L73").Value = 22615.428
$ws.Range("M73").Value = -26202.3
$ws.Range("N73").Value = -24487.428
$ws.Range("H102").Value = 1636.1212
$ws.Range("I102").Value = 1556.3928
$ws.Range("K102").Value = 1556.3928
$ws.Range("M102").Value = 65.60719999999992
$ws.Range("H122").Value = 8427.299999999999
$ws.Range("I122").Value = 9449.875
$ws.Range("J122").Value = 4337
$ws.Range("K122").Value = 28349.625
$ws.Range("L122").Value = 13011
$ws.Range("M122").Value = -25899.625
$ws.Range("N122").Value = -17911
$ws.Range("H126").Value = 11493.625
$ws.Range("J126").Value = 3057
$ws.Range("L126").Value = 9171
$ws.Range("N126").Value = -14111
$ws.Range("H133").Value = 80999.5
$ws.Range("J133").Value = 80999.5
$ws.Range("L133").Value = 80999.5
$ws.Range("N133").Value = -91119.5
$ws.Range("H135").Value = 74661.664
$ws.Range("J135").Value = 74661.664
$ws.Range("L135").Value = 74661.664
$ws.Range("N135").Value = -84801.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10885.286
$ws.Range("I7").Value = 13052.2
$ws.Range("K7").Value = 13052.2
$ws.Range("M7").Value = -12940.2
$ws.Range("H40").Value = 3244.7083
$ws.Range("I40").Value = 2670.2104
$ws.Range("J40").Value = 5427.8
$ws.Range("K40").Value = 2670.2104
$ws.Range("L40").Value = 5427.8
$ws.Range("M40").Value = -2534.2104
$ws.Range("N40").Value = -5699.8
$ws.Range("H126").Value = 10885.286
$ws.Range("I126").Value = 13052.2
$ws.Range("K126").Value = 39156.60000000001
$ws.Range("M126").Value = -36686.60000000001
$ws.Range("H129").Value = 69998.5
$ws.Range("J129").Value = 69998.5
$ws.Range("L129").Value = 69998.5
$ws.Range("N129").Value = -79998.5
$ws.Range("H132").Value = 3520905.2
$ws.Range("I132").Value = 3934635.2
$ws.Range("K132").Value = 11803905.6
$ws.Range("M132").Value = -11801375.6
$ws.Range("H133").Value = 74650.664
$ws.Range("J133").Value = 74650.664
$ws.Range("L133").Value = 74650.664
$ws.Range("N133").Value = -79710.664
$ws.Range("H134").Value = 60420
$ws.Range("J134").Value = 60420
$ws.Range("L134").Value = 60420
$ws.Range("N134").Value = -70560
$ws.Range("H136").Value = 5213787
$ws.Range("I136").Value = 3679543
$ws.Range("J136").Value = 8939808
$ws.Range("K136").Value = 11038629
$ws.Range("L136").Value = 26819424
$ws.Range("M136").Value = -11036079
$ws.Range("N136").Value = -26824524

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3175671.8
$ws.Range("J113").Value = 8548014
$ws.Range("L113").Value = 25644042
$ws.Range("N113").Value = -25648382
$ws.Range("H133").Value = 45357.5
$ws.Range("J133").Value = 45357.5
$ws.Range("L133").Value = 45357.5
$ws.Range("N133").Value = -55477.5
